$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Shared / rich-text strings
#   - "Volume 30   Number  34"  -> "...Number  35"
#   - "Report Covering the Week  8/21/2023  Through  8/27/2023"
#        -> "...8/28/2023  Through  9/3/2023"
# These cells store the text as several rich-text runs that all share the
# exact same font (sz 10, Andale WT). We update only the specific
# sub-strings that changed (via Characters) and then re-apply the matching
# font to the edited run so that the edited fragment keeps explicit
# run-level formatting consistent with the rest of the text.
# ---------------------------------------------------------------------------

$volCell = $ws.Range("A8")
$volChars = $volCell.Characters(21, 2)
$volChars.Text = "35"
$volChars.Font.Name = "Andale WT"
$volChars.Font.Size = 10

$weekCell = $ws.Range("C9")
$firstDate = $weekCell.Characters(27, 9)
$firstDate.Text = "8/28/2023"
$secondDate = $weekCell.Characters(47, 9)
$secondDate.Text = "9/3/2023"
$firstDate = $weekCell.Characters(27, 9)
$firstDate.Font.Name = "Andale WT"
$firstDate.Font.Size = 10
$secondDate = $weekCell.Characters(47, 8)
$secondDate.Font.Name = "Andale WT"
$secondDate.Font.Size = 10

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -76.923076923076

# ---------------------------------------------------------------------------
# Row 15 - Rape
#   C15 switches from a number (1) to the "0" placeholder text used
#   throughout the sheet for empty/zero counts. C14 already holds that
#   exact placeholder (style + shared string), so copy from there.
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("N15").Value = -70.370370370370

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 92.857142857142
$ws.Range("I16").Value = 179
$ws.Range("J16").Value = 161
$ws.Range("K16").Value = 11.180124223602
$ws.Range("L16").Value = 26.056338028169
$ws.Range("M16").Value = 20.134228187919
$ws.Range("N16").Value = -69.349315068493

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -48.648648648648
$ws.Range("I17").Value = 267
$ws.Range("J17").Value = 321
$ws.Range("K17").Value = -16.822429906542
$ws.Range("L17").Value = 12.658227848101
$ws.Range("M17").Value = 64.814814814814
$ws.Range("N17").Value = -41.575492341356

# ---------------------------------------------------------------------------
# Row 18 - Burglary
#   C18 switches from the "0" placeholder text to a real number (2).
#   I14 already holds the plain numeric style used for these counts.
# ---------------------------------------------------------------------------
$ws.Range("I14").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("I18").Value = 80
$ws.Range("J18").Value = 146
$ws.Range("K18").Value = -45.205479452054
$ws.Range("L18").Value = 33.333333333333
$ws.Range("M18").Value = 5.263157894736
$ws.Range("N18").Value = -76.744186046511

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 11.320754716981
$ws.Range("I19").Value = 340
$ws.Range("J19").Value = 367
$ws.Range("K19").Value = -7.356948228882
$ws.Range("L19").Value = 64.251207729468
$ws.Range("M19").Value = 131.292517006803
$ws.Range("N19").Value = 23.636363636363

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
#   C20 switches from a number (2) to the "0" placeholder text.
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 25
$ws.Range("J20").Value = 56
$ws.Range("K20").Value = -1.785714285714
$ws.Range("M20").Value = 12.244897959183
$ws.Range("N20").Value = -75.555555555555

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 2.777777777777
$ws.Range("F21").Value = 118
$ws.Range("G21").Value = 128
$ws.Range("H21").Value = -7.8125
$ws.Range("I21").Value = 935
$ws.Range("J21").Value = 1066
$ws.Range("K21").Value = -12.288930581613
$ws.Range("L21").Value = 31.320224719101
$ws.Range("M21").Value = 55.058043117744
$ws.Range("N21").Value = -51.754385964912

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = -46.666666666666
$ws.Range("L22").Value = -20

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 11
$ws.Range("H23").Value = -26.666666666666
$ws.Range("I23").Value = 113
$ws.Range("J23").Value = 119
$ws.Range("K23").Value = -5.042016806722
$ws.Range("L23").Value = 9.708737864077
$ws.Range("M23").Value = 56.944444444444

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -18.518518518518
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -1.612903225806
$ws.Range("I24").Value = 859
$ws.Range("J24").Value = 885
$ws.Range("K24").Value = -2.937853107344
$ws.Range("L24").Value = 49.651567944250
$ws.Range("M24").Value = 36.565977742448

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -32.352941176470
$ws.Range("I25").Value = 338
$ws.Range("J25").Value = 376
$ws.Range("K25").Value = -10.106382978723
$ws.Range("L25").Value = 4
$ws.Range("M25").Value = -0.879765395894

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
#   C26 switches from a number (1) to the "0" placeholder text.
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("C26"))

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 38
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = -24
$ws.Range("L27").Value = -20.833333333333

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic.
#   C28 switches from the "0" placeholder text to a real number (2).
# ---------------------------------------------------------------------------
$ws.Range("I14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2
$ws.Range("I28").Value = 25
$ws.Range("K28").Value = -21.875
$ws.Range("L28").Value = 25
$ws.Range("M28").Value = -13.793103448275
$ws.Range("N28").Value = -47.916666666666

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc.
#   C29 switches from the "0" placeholder text to a real number (1).
# ---------------------------------------------------------------------------
$ws.Range("I14").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 23
$ws.Range("K29").Value = 4.545454545454
$ws.Range("L29").Value = 35.294117647058
$ws.Range("M29").Value = -4.166666666666
$ws.Range("N29").Value = -51.063829787234
